$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf1"
$ws.Range("C2").Value = "Fgfr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6873773333333334
$ws.Range("H2").Value = 2.062132
$ws.Range("I2").Value = 0.02660947569874856
$ws.Range("J2").Value = 0.02660947569874856
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.132884
$ws.Range("N2").Value = 0.398652
$ws.Range("O2").Value = 0.01195569974366677
$ws.Range("P2").Value = 0.01195569974366677
$ws.Range("Q2").Value = 0.09134144956266668
$ws.Range("R2").Value = 0.822073046064
$ws.Range("S2").Value = 0.0003181349017906352
$ws.Range("T2").Value = 0.0003181349017906352
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf1"
$ws.Range("C3").Value = "Fgfr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6873773333333334
$ws.Range("H3").Value = 2.062132
$ws.Range("I3").Value = 0.02660947569874856
$ws.Range("J3").Value = 0.02660947569874856
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.07352966666666667
$ws.Range("N3").Value = 0.220589
$ws.Range("O3").Value = 0.006615533976389704
$ws.Range("P3").Value = 0.006615533976389703
$ws.Range("Q3").Value = 0.05054262619422223
$ws.Range("R3").Value = 0.4548836357480001
$ws.Range("S3").Value = 0.0001760358905789873
$ws.Range("T3").Value = 0.0001760358905789872
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf1"
$ws.Range("C4").Value = "Fgfr4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6873773333333334
$ws.Range("H4").Value = 2.062132
$ws.Range("I4").Value = 0.02660947569874856
$ws.Range("J4").Value = 0.02660947569874856
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.908285
$ws.Range("N4").Value = 32.724855
$ws.Range("O4").Value = 0.9814287662799436
$ws.Range("P4").Value = 0.9814287662799435
$ws.Range("Q4").Value = 7.498107854540001
$ws.Range("R4").Value = 67.48297069086
$ws.Range("S4").Value = 0.02611530490637894
$ws.Range("T4").Value = 0.02611530490637894
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf1"
$ws.Range("C5").Value = "Fgfr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.913984666666667
$ws.Range("H5").Value = 17.741954
$ws.Range("I5").Value = 0.2289398029860915
$ws.Range("J5").Value = 0.2289398029860915
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.132884
$ws.Range("N5").Value = 0.398652
$ws.Range("O5").Value = 0.01195569974366677
$ws.Range("P5").Value = 0.01195569974366677
$ws.Range("Q5").Value = 0.7858739384453334
$ws.Range("R5").Value = 7.072865446008
$ws.Range("S5").Value = 0.002737135543875933
$ws.Range("T5").Value = 0.002737135543875934
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf1"
$ws.Range("C6").Value = "Fgfr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.913984666666667
$ws.Range("H6").Value = 17.741954
$ws.Range("I6").Value = 0.2289398029860915
$ws.Range("J6").Value = 0.2289398029860915
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.07352966666666667
$ws.Range("N6").Value = 0.220589
$ws.Range("O6").Value = 0.006615533976389704
$ws.Range("P6").Value = 0.006615533976389703
$ws.Range("Q6").Value = 0.4348533212117778
$ws.Range("R6").Value = 3.913679890906
$ws.Range("S6").Value = 0.001514559045202453
$ws.Range("T6").Value = 0.001514559045202453
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf1"
$ws.Range("C7").Value = "Fgfr4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.913984666666667
$ws.Range("H7").Value = 17.741954
$ws.Range("I7").Value = 0.2289398029860915
$ws.Range("J7").Value = 0.2289398029860915
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.908285
$ws.Range("N7").Value = 32.724855
$ws.Range("O7").Value = 0.9814287662799436
$ws.Range("P7").Value = 0.9814287662799435
$ws.Range("Q7").Value = 64.51143022963
$ws.Range("R7").Value = 580.60287206667
$ws.Range("S7").Value = 0.2246881083970131
$ws.Range("T7").Value = 0.2246881083970131
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf1"
$ws.Range("C8").Value = "Fgfr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 19.230689
$ws.Range("H8").Value = 57.692067
$ws.Range("I8").Value = 0.7444507213151601
$ws.Range("J8").Value = 0.7444507213151601
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.132884
$ws.Range("N8").Value = 0.398652
$ws.Range("O8").Value = 0.01195569974366677
$ws.Range("P8").Value = 0.01195569974366677
$ws.Range("Q8").Value = 2.555450877076
$ws.Range("R8").Value = 22.999057893684
$ws.Range("S8").Value = 0.008900429298000199
$ws.Range("T8").Value = 0.008900429298000199
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf1"
$ws.Range("C9").Value = "Fgfr4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 19.230689
$ws.Range("H9").Value = 57.692067
$ws.Range("I9").Value = 0.7444507213151601
$ws.Range("J9").Value = 0.7444507213151601
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.07352966666666667
$ws.Range("N9").Value = 0.220589
$ws.Range("O9").Value = 0.006615533976389704
$ws.Range("P9").Value = 0.006615533976389703
$ws.Range("Q9").Value = 1.414026151940334
$ws.Range("R9").Value = 12.726235367463
$ws.Range("S9").Value = 0.004924939040608265
$ws.Range("T9").Value = 0.004924939040608264
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf1"
$ws.Range("C10").Value = "Fgfr4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 19.230689
$ws.Range("H10").Value = 57.692067
$ws.Range("I10").Value = 0.7444507213151601
$ws.Range("J10").Value = 0.7444507213151601
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.908285
$ws.Range("N10").Value = 32.724855
$ws.Range("O10").Value = 0.9814287662799436
$ws.Range("P10").Value = 0.9814287662799435
$ws.Range("Q10").Value = 209.773836358365
$ws.Range("R10").Value = 1887.964527225285
$ws.Range("S10").Value = 0.7306253529765516
$ws.Range("T10").Value = 0.7306253529765516
